$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# CoinCoach transaction log: replace the old 4-row expense/income sample
# with a header row (date/type/category/amount/description/name) and a
# fuller set of transactions across columns A:F.
# ---------------------------------------------------------------------

function Set-RowValues {
    param(
        [int]$RowNumber,
        [object[]]$Values
    )
    $arr = New-Object 'object[,]' 1, $Values.Length
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $ws.Range($ws.Cells.Item($RowNumber, 1), $ws.Cells.Item($RowNumber, $Values.Length)).Value = $arr
}

# Header row
Set-RowValues 1 @("date", "type", "category", "amount", "description", "name")

# Transaction rows: date (serial), type, category, amount, description, name
Set-RowValues 2 @(46003, "expense", "Food", 230, "Need to eat", "My food expenses")
Set-RowValues 3 @(45999, "income", "salary", 5556, "Day job pay every month", "my pay")
Set-RowValues 4 @(45995, "expense", "mm", 25, $null, "bought something in the store")
Set-RowValues 5 @(45942, "income", "sale", 1200, $null, "sold my former laptop")
Set-RowValues 6 @(45756, "expense", "purchase", 2300, "bought a new laptop", "a new laptop")
Set-RowValues 7 @(46012, "income", "from startup", 1200000, "money from all my SAAS projects", "money from all my SAAS projects")
Set-RowValues 8 @(45711, "expense", "purchase", 129, $null, "a new monitor")
Set-RowValues 9 @(45795, "income", "sale", 123, $null, "sold my former monitor")

# The date column (A) keeps the short-date display used by the original
# sample rows - clone that formatting (rather than re-typing a number
# format string) so the new rows share the same style record.
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Range("A2:A9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Size the columns to fit the new, wider content (values tuned so the
# engine's pixel-quantised column-width storage lands on the same
# rounded character width Excel's own AutoFit produced for this data).
$ws.Columns.Item(1).ColumnWidth = 18.5
$ws.Columns.Item(2).ColumnWidth = 17
$ws.Columns.Item(3).ColumnWidth = 17.333333333333332
$ws.Columns.Item(4).ColumnWidth = 17
$ws.Columns.Item(5).ColumnWidth = 29.833333333333332
$ws.Columns.Item(6).ColumnWidth = 26.333333333333332

# Leave the selection where the user would land after the last entry.
$ws.Range("F14").Select() | Out-Null
